# Auto-generated edit script: scheduled-runner market data refresh
# Updates literal value cells (no formulas) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (44 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1451.76
$ws.Range("I15").Value = 1451.76
$ws.Range("K15").Value = 4355.28
$ws.Range("M15").Value = -4186.28
$ws.Range("H33").Value = 812.1875
$ws.Range("J33").Value = 300
$ws.Range("L33").Value = 300
$ws.Range("N33").Value = -758
$ws.Range("H40").Value = 5590.6
$ws.Range("I40").Value = 2638.1667
$ws.Range("J40").Value = 10019.25
$ws.Range("K40").Value = 2638.1667
$ws.Range("L40").Value = 10019.25
$ws.Range("M40").Value = -2463.1667
$ws.Range("N40").Value = -10369.25
$ws.Range("H62").Value = 2299.1667
$ws.Range("I62").Value = 1759
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1759
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1135
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 2299.1667
$ws.Range("I65").Value = 1759
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 8795
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -5675
$ws.Range("N65").Value = -31240
$ws.Range("H86").Value = 4033.5557
$ws.Range("I86").Value = 3662.75
$ws.Range("K86").Value = 3662.75
$ws.Range("M86").Value = -2539.75
$ws.Range("H89").Value = 4033.5557
$ws.Range("I89").Value = 3662.75
$ws.Range("K89").Value = 18313.75
$ws.Range("M89").Value = -12697.75
$ws.Range("H132").Value = 18548
$ws.Range("I132").Value = 1518.3334
$ws.Range("J132").Value = 86666.664
$ws.Range("K132").Value = 4555.0002
$ws.Range("L132").Value = 259999.992
$ws.Range("M132").Value = -2025.0002
$ws.Range("N132").Value = -265059.992

# --- Sheet: ARM (44 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15809.81
$ws.Range("I2").Value = 16449.65
$ws.Range("J2").Value = 3013
$ws.Range("K2").Value = 16449.65
$ws.Range("L2").Value = 3013
$ws.Range("M2").Value = -16336.65
$ws.Range("N2").Value = -3239
$ws.Range("H32").Value = 8408.556
$ws.Range("I32").Value = 8408.556
$ws.Range("K32").Value = 8408.556
$ws.Range("M32").Value = -8121.556
$ws.Range("H45").Value = 3056.0908
$ws.Range("I45").Value = 2479.7144
$ws.Range("K45").Value = 2479.7144
$ws.Range("M45").Value = -2102.7144
$ws.Range("H61").Value = 4312.1665
$ws.Range("I61").Value = 1681.2222
$ws.Range("J61").Value = 6943.1113
$ws.Range("K61").Value = 1681.2222
$ws.Range("L61").Value = 6943.1113
$ws.Range("M61").Value = -1469.2222
$ws.Range("N61").Value = -7367.1113
$ws.Range("H116").Value = 15809.81
$ws.Range("I116").Value = 16449.65
$ws.Range("J116").Value = 3013
$ws.Range("K116").Value = 16449.65
$ws.Range("L116").Value = 3013
$ws.Range("M116").Value = -14155.65
$ws.Range("N116").Value = -7601
$ws.Range("H119").Value = 89999
$ws.Range("J119").Value = 89999
$ws.Range("L119").Value = 89999
$ws.Range("N119").Value = -99675
$ws.Range("H132").Value = 1046.2778
$ws.Range("I132").Value = 855.6
$ws.Range("K132").Value = 2566.8
$ws.Range("M132").Value = -36.80000000000018
$ws.Range("H136").Value = 4312.1665
$ws.Range("I136").Value = 1681.2222
$ws.Range("J136").Value = 6943.1113
$ws.Range("K136").Value = 5043.6666
$ws.Range("L136").Value = 20829.3339
$ws.Range("M136").Value = -2493.6666
$ws.Range("N136").Value = -25929.3339

# --- Sheet: BSM (29 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15809.81
$ws.Range("I3").Value = 16449.65
$ws.Range("J3").Value = 3013
$ws.Range("K3").Value = 16449.65
$ws.Range("L3").Value = 3013
$ws.Range("M3").Value = -16335.65
$ws.Range("N3").Value = -3241
$ws.Range("H86").Value = 3625.0386
$ws.Range("I86").Value = 1647.55
$ws.Range("J86").Value = 10216.667
$ws.Range("K86").Value = 1647.55
$ws.Range("L86").Value = 10216.667
$ws.Range("M86").Value = -524.55
$ws.Range("N86").Value = -12462.667
$ws.Range("H89").Value = 3625.0386
$ws.Range("I89").Value = 1647.55
$ws.Range("J89").Value = 10216.667
$ws.Range("K89").Value = 8237.75
$ws.Range("L89").Value = 51083.335
$ws.Range("M89").Value = -2621.75
$ws.Range("N89").Value = -62315.335
$ws.Range("H94").Value = 3605.9167
$ws.Range("I94").Value = 2717.2
$ws.Range("K94").Value = 2717.2
$ws.Range("M94").Value = -2266.2
$ws.Range("H134").Value = 3709.5625
$ws.Range("I134").Value = 2248.5652
$ws.Range("K134").Value = 6745.6956
$ws.Range("M134").Value = -4210.6956

# --- Sheet: CRP (29 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4132.636
$ws.Range("I31").Value = 1201.5
$ws.Range("J31").Value = 7650
$ws.Range("K31").Value = 1201.5
$ws.Range("L31").Value = 7650
$ws.Range("M31").Value = -906.5
$ws.Range("N31").Value = -8240
$ws.Range("H34").Value = 4132.636
$ws.Range("I34").Value = 1201.5
$ws.Range("J34").Value = 7650
$ws.Range("K34").Value = 1201.5
$ws.Range("L34").Value = 7650
$ws.Range("M34").Value = -999.5
$ws.Range("N34").Value = -8054
$ws.Range("H58").Value = 1510.44
$ws.Range("I58").Value = 704.6667
$ws.Range("K58").Value = 704.6667
$ws.Range("M58").Value = -501.6667
$ws.Range("H132").Value = 1869.381
$ws.Range("I132").Value = 1558.7222
$ws.Range("J132").Value = 3733.3333
$ws.Range("K132").Value = 4676.1666
$ws.Range("L132").Value = 11199.9999
$ws.Range("M132").Value = -2146.1666
$ws.Range("N132").Value = -16259.9999
$ws.Range("H136").Value = 1510.44
$ws.Range("I136").Value = 704.6667
$ws.Range("K136").Value = 2114.0001
$ws.Range("M136").Value = 435.9998999999998

# --- Sheet: CUL (11 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2192.75
$ws.Range("I129").Value = 1221.5
$ws.Range("J129").Value = 2775.5
$ws.Range("K129").Value = 3664.5
$ws.Range("L129").Value = 8326.5
$ws.Range("M129").Value = 1335.5
$ws.Range("N129").Value = -18326.5
$ws.Range("H131").Value = 6652.6
$ws.Range("J131").Value = 9199.4
$ws.Range("L131").Value = 27598.2
$ws.Range("N131").Value = -37678.2

# --- Sheet: GSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1889846.2
$ws.Range("J11").Value = 1592857.1
$ws.Range("L11").Value = 1592857.1
$ws.Range("N11").Value = -1593135.1
$ws.Range("H122").Value = 5173.8
$ws.Range("I122").Value = 1316.909
$ws.Range("K122").Value = 3950.727
$ws.Range("M122").Value = -1500.727
$ws.Range("H132").Value = 1858.3334
$ws.Range("I132").Value = 1858.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5575.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3045.0002
$ws.Range("N132").ClearContents()

# --- Sheet: LTW (4 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 973.1429000000001
$ws.Range("I46").Value = 882.6
$ws.Range("K46").Value = 882.6
$ws.Range("M46").Value = -694.6

# --- Sheet: WVR (19 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 20000
$ws.Range("J9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("N9").Value = -20280
$ws.Range("H81").Value = 6590.65
$ws.Range("J81").Value = 1550
$ws.Range("L81").Value = 3100
$ws.Range("N81").Value = -5222
$ws.Range("H84").Value = 6590.65
$ws.Range("J84").Value = 1550
$ws.Range("L84").Value = 15500
$ws.Range("N84").Value = -26108
$ws.Range("H136").Value = 2595.0688
$ws.Range("I136").Value = 1150.2941
$ws.Range("J136").Value = 4641.8335
$ws.Range("K136").Value = 3450.8823
$ws.Range("L136").Value = 13925.5005
$ws.Range("M136").Value = -900.8823000000002
$ws.Range("N136").Value = -19025.5005

Write-Host "Applied 195 cell updates across 8 sheets."